# Cite Osborn 1899 on rearing Diplodocus
#
# Applies two changes:
#  1. Expands the "However, the notion of rearing sauropods..." sentence in
#     the AMNH-mount section to quote Osborn (1899:213) at length, and moves
#     the Hatcher (1901:57-58) citation to follow it.
#  2. Adds a new "Osborn, Henry F. 1899..." entry to the reference list,
#     right after the McIntosh (2005) entry.

function Apply-FormattedRuns($d, $targetRange, $segments) {
    # Build the plain-text concatenation of every segment, drop it into
    # $targetRange (replacing whatever text used to be there), and then walk
    # back over the same character offsets applying per-segment character
    # formatting via small sub-ranges. Clearing the range first and using
    # InsertAfter (rather than assigning .Text on a range that still spans
    # old content) keeps the existing/empty run properties (e.g. an empty
    # <w:rPr/>) intact instead of the engine fabricating bare runs.
    $full = ""
    foreach ($seg in $segments) {
        $full = $full + $seg.text
    }

    $targetRange.Text = ""
    $start = $targetRange.Start
    $targetRange.InsertAfter($full)

    $pos = $start
    foreach ($seg in $segments) {
        $len = $seg.text.Length
        if ($len -gt 0 -and $seg.fmt -ne "plain") {
            $sub = $d.Range($pos, $pos + $len)
            if ($seg.fmt -eq "i") {
                $sub.Font.Italic = 1
            } elseif ($seg.fmt -eq "b") {
                $sub.Font.Bold = 1
            }
        }
        $pos = $pos + $len
    }
    return $pos
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Rewrite the "heritage going back..." sentence.
# ---------------------------------------------------------------------

$range = $d.Content
$oldSentence = " the notion of rearing sauropods has a heritage going back at least to the Hatcher (1901:57–58), who strongly implied that without quite explicitly stating that "
$found = $range.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the original 'heritage going back' sentence"
}

$newSegments = @(
    @{text=" the notion of rearing sauropods has a heritage going back at least to "; fmt="plain"},
    @{text="Osborn (1899:213), who wrote that the tail of "; fmt="plain"},
    @{text="Diplodocus"; fmt="i"},
    @{text=" “functioned as a lever to balance the weight "; fmt="plain"},
    @{text="of the dorsals, anterior limbs, neck, and head, and to raise the entire forward portion of the body upwards. "; fmt="plain"},
    @{text="[…]"; fmt="plain"},
    @{text=" Thus the quadrupedal Dinosaurs occasionally assumed the position characteristic of the bipedal Dinosaurs — namely, a tripodal position, the body supported upon the hind feet and the tail”. "; fmt="plain"},
    @{text="In his classic monograph of "; fmt="plain"},
    @{text="Diplodocus carnegii"; fmt="i"},
    @{text=", "; fmt="plain"},
    @{text="Hatcher (1901:57–58) strongly implied, without quite explicitly stating, that "; fmt="plain"}
)

Apply-FormattedRuns $d $range $newSegments | Out-Null

# ---------------------------------------------------------------------
# 2. Insert the new Osborn (1899) reference-list entry after McIntosh (2005).
# ---------------------------------------------------------------------

$refRange = $d.Content
$foundRef = $refRange.Find.Execute("Indiana University Press, Bloomington, Indiana. 495 pp.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundRef) {
    throw "Could not find the McIntosh (2005) reference paragraph"
}

# Collapse to the end of that paragraph's text (just before its paragraph
# mark), split a new paragraph in after it, then step past the paragraph
# mark to land inside the freshly created (empty) paragraph, which inherits
# the "Reference" style from its neighbours.
$refRange.Collapse(0)
$refRange.InsertParagraphAfter()
$newParaStart = $refRange.End + 1
$newParaRange = $d.Range($newParaStart, $newParaStart)

$refSegments = @(
    @{text="Osborn, H"; fmt="plain"},
    @{text="enry"; fmt="plain"},
    @{text=". F. 1899. A skeleton of "; fmt="plain"},
    @{text="Diplodocus"; fmt="i"},
    @{text=". "; fmt="plain"},
    @{text="Memoirs of the American Museum of Natural History"; fmt="i"},
    @{text=", "; fmt="plain"},
    @{text="1"; fmt="b"},
    @{text=":189–214 "; fmt="plain"},
    @{text="and "; fmt="plain"},
    @{text="plates 24–28."; fmt="plain"}
)

Apply-FormattedRuns $d $newParaRange $refSegments | Out-Null

Write-Host "Edit complete"
